$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the "VERSION 4" banner (row 14, merged A14:M14) to "VERSION 4-6"
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "VERSION 4-6"

# ---------------------------------------------------------------------
# 2) New bolt / nut size tables in I21:M24
# ---------------------------------------------------------------------

# --- Row 21: bolt table header (bold + border, like A2:D2) ---
$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("I21:L21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I21").Value = "Bolt size"
$ws.Range("J21").Value = "Size"
$ws.Range("K21").Value = "Length S"
$ws.Range("L21").Value = "Length L"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("M21").PasteSpecial(-4122) | Out-Null
$ws.Range("M21").Value = "Unit"

# --- Row 22: bolt table data (normal + border, like A3:D3) ---
$ws.Range("A3:D3").Copy() | Out-Null
$ws.Range("I22:L22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "M3"
$ws.Range("J22").HorizontalAlignment = -4152   # xlRight
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = 20
$ws.Range("D3").Copy() | Out-Null
$ws.Range("M22").PasteSpecial(-4122) | Out-Null
$ws.Range("M22").Value = "mm"

# --- Row 23: nut table header (bold + border) ---
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("I23:K23").PasteSpecial(-4122) | Out-Null
$ws.Range("I23").Value = "Nut size"
$ws.Range("J23").Value = "Borders"
$ws.Range("K23").Value = "Height"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("L23").PasteSpecial(-4122) | Out-Null
$ws.Range("L23").Interior.Color = 16777215   # white fill (distinguishes this xf from the plain bold+border one)
$ws.Range("L23").Value = "Width"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("M23").PasteSpecial(-4122) | Out-Null
$ws.Range("M23").Value = "Unit"

# --- Row 24: nut table data (normal + border) ---
$ws.Range("A3:D3").Copy() | Out-Null
$ws.Range("I24:L24").PasteSpecial(-4122) | Out-Null
$ws.Range("I24").Value = "-"
$ws.Range("J24").Value = 6
$ws.Range("K24").Value = 2.4
$ws.Range("L24").Value = 5.4
$ws.Range("D3").Copy() | Out-Null
$ws.Range("M24").PasteSpecial(-4122) | Out-Null
$ws.Range("M24").Value = "mm"

# ---------------------------------------------------------------------
# 3) Columns J:L grew wider to fit the new content
# ---------------------------------------------------------------------
$ws.Columns.Item(10).AutoFit() | Out-Null
$ws.Columns.Item(11).AutoFit() | Out-Null
$ws.Columns.Item(12).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# 4) Selection left on Q27 (as in the saved file)
# ---------------------------------------------------------------------
$ws.Range("Q27").Select() | Out-Null

Write-Host "done"
